$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financements")

# Remove the useless sample line (row 3) from the template, keeping the
# cell formatting/styles intact while clearing the values/content.
$ws.Range("A3:F3").ClearContents()

# Update the selection to match the author's saved state: the whole
# A3:F3 row is now selected (previously only F3 was selected).
$ws.Range("A3:F3").Select()
